# Refresh the cryptocurrency price/volume snapshot (and fix two swapped
# row pairs) to match the "Updated cryptos list" GitHub Actions commit.
#
# Prices in column D are entered with a leading apostrophe when they look
# like numbers (e.g. "1.15", "0.0950") so Excel stores them as literal text
# instead of silently re-parsing/rounding them as numeric values - matching
# the original inline-string cells. The apostrophe adds a quote-prefix
# style to the cell, so afterwards we reset those cells back to the
# "Normal" style to keep the same (unstyled) formatting as before.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.750.09"
$ws.Range("E2").Value = "  -1.04%  "
$ws.Range("D3").Value = "2.075.66"
$ws.Range("E3").Value = "  +1.23%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "'244.69"
$ws.Range("E5").Value = "  -1.38%  "
$ws.Range("D6").Value = "'0.651"
$ws.Range("E6").Value = "  -2.01%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").Value = "'54.93"
$ws.Range("E8").Value = "  -5.64%  "
$ws.Range("D9").Value = "'59.71"
$ws.Range("E9").Value = "  -0.29%  "
$ws.Range("D10").Value = "'0.366"
$ws.Range("E10").Value = "  -4.23%  "
$ws.Range("D11").Value = "'0.0759"
$ws.Range("E11").Value = "  -2.25%  "
$ws.Range("E12").Value = "  +1.12%  "
$ws.Range("D13").Value = "'14.94"
$ws.Range("E13").Value = "  -6.26%  "
$ws.Range("D14").Value = "'0.882"
$ws.Range("E14").Value = "  +3.82%  "
$ws.Range("D15").Value = "2.389.84"
$ws.Range("B16").Value = "Polkadot"
$ws.Range("C16").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D16").Value = "'5.48"
$ws.Range("E16").Value = "  -4.01%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "2.056.94"
$ws.Range("E17").Value = "  +0.31%  "
$ws.Range("D18").Value = "36.717.72"
$ws.Range("E18").Value = "  -1.14%  "
$ws.Range("D19").Value = "'17.26"
$ws.Range("E19").Value = "  -3.97%  "
$ws.Range("D20").Value = "'72.70"
$ws.Range("E20").Value = "  -2.99%  "
$ws.Range("D21").Value = "0.0₃0877"
$ws.Range("E21").Value = "  -1.60%  "
$ws.Range("D22").Value = "'5.41"
$ws.Range("E22").Value = "  +1.19%  "
$ws.Range("D23").Value = "'237.25"
$ws.Range("E23").Value = "  +0.02%  "
$ws.Range("D24").Value = "'0.999"
$ws.Range("E24").Value = "  -0.11%  "
$ws.Range("D25").Value = "'2.40"
$ws.Range("E25").Value = "  -2.30%  "
$ws.Range("D26").Value = "'9.80"
$ws.Range("E26").Value = "  +3.26%  "
$ws.Range("D27").Value = "'2.16"
$ws.Range("E27").Value = "  -0.74%  "
$ws.Range("D28").Value = "'167.03"
$ws.Range("E28").Value = "  -1.46%  "
$ws.Range("D29").Value = "'20.53"
$ws.Range("E29").Value = "  +2.32%  "
$ws.Range("B30").Value = "Stellar"
$ws.Range("C30").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D30").Value = "'0.123"
$ws.Range("E30").Value = "  -1.31%  "
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").Value = "'5.26"
$ws.Range("E31").Value = "  +9.50%  "
$ws.Range("D32").Value = "'1.18"
$ws.Range("E32").Value = "  +6.38%  "
$ws.Range("D33").Value = "'4.66"
$ws.Range("E33").Value = "  +4.21%  "
$ws.Range("D34").Value = "'0.0607"
$ws.Range("E34").Value = "  -1.71%  "
$ws.Range("D35").Value = "'2.36"
$ws.Range("E35").Value = "  +4.93%  "
$ws.Range("E36").Value = "  +0.35%  "
$ws.Range("E37").Value = "  +4.32%  "
$ws.Range("D38").Value = "'0.0833"
$ws.Range("E38").Value = "  -6.89%  "
$ws.Range("D39").Value = "'1.27"
$ws.Range("E39").Value = "  -4.53%  "
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").Value = "'0.0220"
$ws.Range("E40").Value = "  -0.93%  "
$ws.Range("B41").Value = "ARBITRUM"
$ws.Range("C41").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D41").Value = "'1.15"
$ws.Range("E41").Value = "  +0.88%  "
$ws.Range("D42").Value = "'4.85"
$ws.Range("E42").Value = "  -6.51%  "
$ws.Range("D43").Value = "'0.0950"
$ws.Range("E43").Value = "  -4.33%  "
$ws.Range("D44").Value = "'96.24"
$ws.Range("E44").Value = "  +0.51%  "
$ws.Range("D45").Value = "'2.85"
$ws.Range("E45").Value = "  -12.91%  "
$ws.Range("D46").Value = "'15.99"
$ws.Range("E46").Value = "  -7.42%  "
$ws.Range("D47").Value = "1.348.81"
$ws.Range("E47").Value = "  +5.70%  "
$ws.Range("D48").Value = "'2.42"
$ws.Range("E48").Value = "  -1.12%  "
$ws.Range("D49").Value = "'7.21"
$ws.Range("E49").Value = "  +5.99%  "
$ws.Range("D50").Value = "'2.89"
$ws.Range("E50").Value = "  +1.44%  "
$ws.Range("D51").Value = "2.269.63"
$ws.Range("E51").Value = "  +1.77%  "

# Clear the quote-prefix style added by forcing numeric-looking strings as text,
# restoring the default (unstyled) cell formatting to match the original cells.
$textRefs = @("D5","D6","D8","D9","D10","D11","D13","D14","D16","D19","D20","D22","D23","D24","D25","D26","D27","D28","D29","D30","D31","D32","D33","D34","D35","D38","D39","D40","D41","D42","D43","D44","D45","D46","D48","D49","D50")
foreach ($r in $textRefs) {
    $ws.Range($r).Style = "Normal"
}
